$wb = $excel.ActiveWorkbook


# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H8").Value = 131.875
$ws.Range("I8").Value = 169.16667
$ws.Range("J8").Value = 20
$ws.Range("K8").Value = 507.50001
$ws.Range("L8").Value = 60
$ws.Range("M8").Value = -368.50001
$ws.Range("N8").Value = -338
$ws.Range("H44").Value = 12100
$ws.Range("J44").Value = 12100
$ws.Range("L44").Value = 12100
$ws.Range("N44").Value = -13024
$ws.Range("H118").Value = 595.8
$ws.Range("I118").Value = 276.33334
$ws.Range("J118").Value = 1075
$ws.Range("K118").Value = 829.0000200000001
$ws.Range("L118").Value = 3225
$ws.Range("M118").Value = 827.9999799999999
$ws.Range("N118").Value = -6539
$ws.Range("H123").Value = 29250.857
$ws.Range("J123").Value = 29250.857
$ws.Range("L123").Value = 29250.857
$ws.Range("N123").Value = -39050.857

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H19").Value = 1000
$ws.Range("I19").Value = 1000
$ws.Range("J19").Value = 0
$ws.Range("K19").Value = 1000
$ws.Range("L19").Value = 0
$ws.Range("M19").Value = -771
$ws.Range("N19").ClearContents()
$ws.Range("H32").Value = 3857.4
$ws.Range("I32").Value = 2626.0186
$ws.Range("J32").Value = 8013.3125
$ws.Range("K32").Value = 2626.0186
$ws.Range("L32").Value = 8013.3125
$ws.Range("M32").Value = -2339.0186
$ws.Range("N32").Value = -8587.3125
$ws.Range("H86").Value = 43647.332
$ws.Range("J86").Value = 43647.332
$ws.Range("L86").Value = 43647.332
$ws.Range("N86").Value = -46019.332
$ws.Range("H89").Value = 43647.332
$ws.Range("J89").Value = 43647.332
$ws.Range("L89").Value = 130941.996
$ws.Range("N89").Value = -142797.996
$ws.Range("H101").Value = 49500
$ws.Range("J101").Value = 49500
$ws.Range("L101").Value = 49500
$ws.Range("N101").Value = -55990
$ws.Range("H132").Value = 9061.227999999999
$ws.Range("I132").Value = 2150.6667
$ws.Range("J132").Value = 17353.9
$ws.Range("K132").Value = 6452.000100000001
$ws.Range("L132").Value = 52061.7
$ws.Range("M132").Value = -3922.000100000001
$ws.Range("N132").Value = -57121.7

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H18").Value = 0
$ws.Range("J18").Value = 0
$ws.Range("L18").Value = 0
$ws.Range("N18").ClearContents()
$ws.Range("H99").Value = 1000000000
$ws.Range("I99").Value = 1000000000
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 1000000000
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = -999998502
$ws.Range("N99").ClearContents()

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H54").Value = 8000
$ws.Range("J54").Value = 8000
$ws.Range("L54").Value = 8000
$ws.Range("N54").Value = -9316
$ws.Range("H116").Value = 40720
$ws.Range("J116").Value = 40720
$ws.Range("L116").Value = 40720
$ws.Range("N116").Value = -49898
$ws.Range("H118").Value = 39690.43
$ws.Range("J118").Value = 39690.43
$ws.Range("L118").Value = 39690.43
$ws.Range("N118").Value = -43004.43
$ws.Range("H119").Value = 50000
$ws.Range("J119").Value = 50000
$ws.Range("L119").Value = 50000
$ws.Range("N119").Value = -59676
$ws.Range("H120").Value = 0
$ws.Range("J120").Value = 0
$ws.Range("L120").Value = 0
$ws.Range("N120").ClearContents()
$ws.Range("H121").Value = 46330
$ws.Range("J121").Value = 46330
$ws.Range("L121").Value = 46330
$ws.Range("N121").Value = -48950

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H96").Value = 6420
$ws.Range("J96").Value = 6420
$ws.Range("L96").Value = 19260
$ws.Range("N96").Value = -23378
$ws.Range("H125").Value = 2995.7273
$ws.Range("J125").Value = 3145.3
$ws.Range("L125").Value = 9435.900000000001
$ws.Range("N125").Value = -19275.9
$ws.Range("H131").Value = 44000836
$ws.Range("J131").Value = 66667732
$ws.Range("L131").Value = 200003196
$ws.Range("N131").Value = -200013276
$ws.Range("H132").Value = 1870.7354
$ws.Range("I132").Value = 1884.8
$ws.Range("J132").Value = 1868.3103
$ws.Range("K132").Value = 16963.2
$ws.Range("L132").Value = 16814.7927
$ws.Range("M132").Value = -14433.2
$ws.Range("N132").Value = -21874.7927

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 8132.0625
$ws.Range("I126").Value = 10051
$ws.Range("J126").Value = 2375.25
$ws.Range("K126").Value = 30153
$ws.Range("L126").Value = 7125.75
$ws.Range("M126").Value = -27683
$ws.Range("N126").Value = -12065.75
$ws.Range("H132").Value = 3750.5293
$ws.Range("I132").Value = 2703.6667
$ws.Range("J132").Value = 3974.8572
$ws.Range("K132").Value = 8111.000100000001
$ws.Range("L132").Value = 11924.5716
$ws.Range("M132").Value = -5581.000100000001
$ws.Range("N132").Value = -16984.5716

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 4117644.2
$ws.Range("I22").Value = 27778198
$ws.Range("J22").Value = 2765.1738
$ws.Range("K22").Value = 27778198
$ws.Range("L22").Value = 2765.1738
$ws.Range("M22").Value = -27777903
$ws.Range("N22").Value = -3355.1738
$ws.Range("H27").Value = 4117644.2
$ws.Range("I27").Value = 27778198
$ws.Range("J27").Value = 2765.1738
$ws.Range("K27").Value = 27778198
$ws.Range("L27").Value = 2765.1738
$ws.Range("M27").Value = -27778091
$ws.Range("N27").Value = -2979.1738
$ws.Range("H46").Value = 19608978
$ws.Range("I46").Value = 27778690
$ws.Range("J46").Value = 1669.6
$ws.Range("K46").Value = 27778690
$ws.Range("L46").Value = 1669.6
$ws.Range("M46").Value = -27778502
$ws.Range("N46").Value = -2045.6
$ws.Range("H55").Value = 60000370
$ws.Range("I55").Value = 111111336
$ws.Range("J55").Value = 31250446
$ws.Range("K55").Value = 111111336
$ws.Range("L55").Value = 31250446
$ws.Range("M55").Value = -111111163
$ws.Range("N55").Value = -31250792

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H12").Value = 10000
$ws.Range("J12").Value = 10000
$ws.Range("L12").Value = 10000
$ws.Range("N12").Value = -10284
$ws.Range("H54").Value = 7003.5293
$ws.Range("J54").Value = 7003.5293
$ws.Range("L54").Value = 7003.5293
$ws.Range("N54").Value = -8043.5293